# Fix spelling/accents on slides 4 and 5 (Analisis -> Análisis, Reflexion... -> Reflexión...,
# Presentacion -> Presentación), merging the two runs in the "Reflexion" paragraph into one.

$p = $ppt.ActivePresentation

# Slide 4: "Analisis" -> "Análisis" (title placeholder)
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Análisis"

# Slide 4: content placeholder, first paragraph: "Reflexion" + " sobre los resultados encontrados"
# becomes a single run "Reflexión sobre los resultados encontrados"
$tr = $s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1)
$tr.Text = "Reflexión sobre los resultados encontrados"

# Slide 5: "Presentacion" -> "Presentación" (title placeholder)
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Presentación"
